# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates column G ("K") values for rows 2-19 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    3  = 1
    4  = 2
    5  = 2
    6  = 2
    7  = 1
    8  = 0
    9  = 0
    10 = 2
    11 = 1
    12 = 1
    13 = 2
    14 = 3
    15 = 2
    16 = 3
    17 = 1
    18 = 4
    19 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
